# "Generate Report for Handback" - refresh the handback status report with
# the latest handoff/handback timestamps for the 4bfa4aa4-... file, for both
# the zh-cn and de-de locales, then re-stamp the Overview sheet's "Latest HO
# Xliff Generate Date" cells (same text, kept in sync with the per-locale
# sheets).

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet: row 2 is the 4bfa4aa4-... source file --------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Cells.Item(2, 8).Value  = "2016-08-27 08:48:04"   # H2 Correspond Handoff Datetime
$wsZhCn.Cells.Item(2, 11).Value = "2016-08-27 08:48:29"   # K2 Correspond Handback DateTime

# --- de-de sheet: row 2 is the 4bfa4aa4-... source file --------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Cells.Item(2, 8).Value  = "2016-08-27 08:48:10"   # H2 Correspond Handoff Datetime
$wsDeDe.Cells.Item(2, 11).Value = "2016-08-27 08:48:35"   # K2 Correspond Handback DateTime

# --- Overview sheet: re-write the "Latest HO Xliff Generate Date" column ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Cells.Item(2, 7).Value = "2016-08-27 08:47:13"   # G2 (4bfa4aa4 row)
$wsOverview.Cells.Item(3, 7).Value = "2016-08-27 08:47:13"   # G3 (c775bb97 row)
